# Insert one new data row right before the current row 449 (pushing the
# existing rows 449-508 down to 450-509) and populate the new row with the
# values for the new "Primera" Repollo record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(449).Insert()

$ws.Range("A449").Value2 = 5
$ws.Range("B449").Value2 = "Macroferia Regional de Talca"
$ws.Range("C449").Value2 = "Maule"
$ws.Range("D449").Value2 = 45124
$ws.Range("E449").Value2 = 7
$ws.Range("F449").Value2 = 100112006
$ws.Range("G449").Value2 = "Repollo"
$ws.Range("H449").Value2 = "Crespo record"
$ws.Range("I449").Value2 = "Primera"
$ws.Range("J449").Value2 = 5000
$ws.Range("K449").Value2 = 600
$ws.Range("L449").Value2 = 600
$ws.Range("M449").Value2 = 600
$ws.Range("N449").Value2 = "`$/unidad"
$ws.Range("O449").Value2 = "Región del Maule"
$ws.Range("P449").Value2 = 600
$ws.Range("Q449").Value2 = 1
$ws.Range("R449").Value2 = "Hortaliza"
